$d = $word.ActiveDocument

# Change 1: merge "During the process of dev" + "eloping " into one run,
# removing the _GoBack bookmark that sat between them.
$d.Content.Find.Execute("During the process of dev" + [char]0x0B + "eloping ", $true, $false, $false, $false, $false, $true, 1, $false, "During the process of developing ", 2) | Out-Null

# Change 3: "scheme" -> "schema" and "functional of" -> "functionality of"
$d.Content.Find.Execute("management scheme and functional of", $true, $false, $false, $false, $false, $true, 1, $false, "management schema and functionality of", 2) | Out-Null
